# gh-pages data refresh: the now-stale lead exhibition row
# ("上饶·囧喵喵次元国风动漫游戏展") drops out of the "展览" and "全部类型"
# listings. Every later row moves up one slot to fill the gap (the row's
# own running-index in column A is left as-is — it is positional, not
# carried with the row's content), the last row disappears, and a
# handful of "想去人数" (want-to-go) counts are refreshed upward.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# New-row-number (after the shift) -> refreshed "想去人数" (column F) value.
$fUpdates = @{
    4  = 4703
    6  = 417
    7  = 1421
    8  = 930
    9  = 59
    10 = 1247
    12 = 859
    14 = 68
    15 = 32
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Pull columns B..I up by one row (row 3 -> row 2, row 4 -> row 3, ...,
    # row 18 -> row 17). Column A is left untouched since it is just the
    # 0-based running index for the row's position, not event data.
    # Walking top-down is safe: row N only ever reads from row N+1, which
    # hasn't been overwritten yet at that point in the loop.
    for ($r = 2; $r -le 17; $r++) {
        $srcRow = $r + 1
        for ($c = 2; $c -le 9; $c++) {
            $srcCell = $ws.Cells.Item($srcRow, $c)
            $dstCell = $ws.Cells.Item($r, $c)
            $val = $srcCell.Value2
            if ($val -is [string]) {
                # Route the write through a Text-formatted cell so Excel
                # doesn't auto-coerce date-looking strings (e.g.
                # "2024.02.24") into date serial numbers, then drop the
                # temporary formatting so the cell stays styled exactly
                # like its un-styled siblings.
                $dstCell.NumberFormat = "@"
                $dstCell.Value2 = $val
                $dstCell.ClearFormats()
            } else {
                $dstCell.Value2 = $val
            }
        }
    }

    # The old last row (18) is now a duplicate of row 17; remove it so the
    # sheet's used range shrinks back to A1:I17.
    $ws.Rows.Item(18).Delete()

    # Refresh the "想去人数" counts that changed since the last snapshot.
    foreach ($rowNum in $fUpdates.Keys) {
        $ws.Cells.Item($rowNum, 6).Value2 = $fUpdates[$rowNum]
    }
}
